# The sheet lists English phrases in column A, one per row, with a
# "Phrase" header in A1. "Hello" sits in A2. This change inserts 20 new
# common greeting phrases right after "Hello" (pushing every subsequent
# phrase down by 20 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 20 blank rows right after row 2 ("Hello"), shifting everything
# below down.
$ws.Rows("3:22").Insert()

# New greeting phrases to drop into the freshly inserted rows.
$greetings = @(
    "Hi",
    "Hey",
    "Hello there",
    "Howdy",
    "Greetings",
    "Good day",
    "Salutations",
    "Hey there",
    "What's up",
    "How's it going",
    "How are you doing",
    "Yo",
    "Hiya",
    "What's happening",
    "Good morning",
    "Good afternoon",
    "Good evening",
    "How's everything",
    "Nice to see you",
    "How do you do"
)

for ($i = 0; $i -lt $greetings.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $greetings[$i]
}

# Carry A2's formatting onto B2 (matches the stray formatted-but-empty
# cell that appears next to "Hello" after the edit).
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
